$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column (C) for rows 2-9 from 2023-09-01 (45170) to 2023-09-05 (45174)
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45174
}
